$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = "LK644532"
$ws.Range("C5").Value = "SCL ENTERPRISES LAUNDRY"
$ws.Range("E5").Value = 1720
$ws.Range("F5").Value = "T"
$ws.Range("H5").Value = 45244.04213364583
$ws.Range("J5").Value = "10/23/23 21:05"
$ws.Range("K5").Value = "10/23/23 21:05"
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = "$1,720 as of 10/23/2023 7:05:45 PM"
$ws.Range("N5").Value = 1660
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0

# Row 7
$ws.Range("H7").ClearContents()
$ws.Range("A7").Value = "L678988"
$ws.Range("C7").Value = "PAYELESS MARKET"
$ws.Range("E7").Value = 2400
$ws.Range("F7").Value = "T"
$ws.Range("I7").Value = "ATM Inactive greater than 2000 minutes"
$ws.Range("J7").Value = "07/20/23 20:09"
$ws.Range("K7").Value = "07/20/23 20:09"
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = "$2,400 as of 7/20/2023 6:09:40 PM"
$ws.Range("N7").Value = 2500
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0

# Row 8
$ws.Range("A8").Value = "LK561655"
$ws.Range("C8").Value = "CRENSHAW CRAVOR #2"
$ws.Range("E8").Value = 2780
$ws.Range("F8").Value = "T"
$ws.Range("I8").Value = "ATM Inactive greater than 48 minutes"
$ws.Range("J8").Value = "01/23/20 08:24"
$ws.Range("K8").Value = "01/23/20 08:24"
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = "$2,780 as of 1/23/2020 6:24:32 AM"
$ws.Range("N8").Value = 2800
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0

# Row 9
$ws.Range("A9").Value = "L474792"
$ws.Range("C9").Value = "NICK SHELL SERVICE"
$ws.Range("E9").Value = 2860
$ws.Range("F9").Value = "T"
$ws.Range("H9").Value = 45243.04213364583
$ws.Range("J9").Value = "10/23/23 13:19"
$ws.Range("K9").Value = "10/23/23 13:19"
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = "$2,860 as of 10/23/2023 11:19:13 AM"
$ws.Range("N9").Value = 2860
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0

# Row 10
$ws.Range("I10").ClearContents()
$ws.Range("A10").Value = "L662336"
$ws.Range("C10").Value = "SB#4 MONA MARKET"
$ws.Range("E10").Value = 3120
$ws.Range("F10").Value = "T"
$ws.Range("H10").Value = 45239.04213364583
$ws.Range("J10").Value = "10/23/23 16:57"
$ws.Range("K10").Value = "10/23/23 16:57"
$ws.Range("L10").Value = 120
$ws.Range("M10").Value = "$3,120 as of 10/23/2023 2:57:12 PM"
$ws.Range("N10").Value = 3120
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0

# Row 11
$ws.Range("A11").Value = "L488595"
$ws.Range("C11").Value = "N S MART"
$ws.Range("E11").Value = 3460
$ws.Range("F11").Value = "T"
$ws.Range("H11").Value = 45263.04213364583
$ws.Range("I11").Value = "ATM Inactive greater than 2000 minutes"
$ws.Range("J11").Value = "10/22/23 16:35"
$ws.Range("K11").Value = "10/22/23 16:35"
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = "$3,460 as of 10/22/2023 2:35:21 PM"
$ws.Range("N11").Value = 3440
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0

# Row 13
$ws.Range("A13").Value = "L688961"
$ws.Range("C13").Value = "MONA MART"
$ws.Range("E13").Value = 4000
$ws.Range("F13").Value = "T"
$ws.Range("H13").Value = 45456.04213364583
$ws.Range("I13").Value = "ATM Inactive greater than 2000 minutes"
$ws.Range("J13").Value = "10/17/23 13:26"
$ws.Range("K13").Value = "10/17/23 13:00"
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = "$4,000 as of 10/17/2023 11:00:09 AM"
$ws.Range("N13").Value = 4000
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0

# Row 14
$ws.Range("A14").Value = "L476340"
$ws.Range("C14").Value = "DONUT & SANDWICH"
$ws.Range("E14").Value = 4040
$ws.Range("F14").Value = "T"
$ws.Range("H14").Value = 45242.04213364583
$ws.Range("J14").Value = "10/24/23 14:06"
$ws.Range("K14").Value = "10/24/23 14:06"
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = "$4,060 as of 10/24/2023 10:00:58 AM"
$ws.Range("N14").Value = 3880
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0

# Row 15
$ws.Range("A15").Value = "L697590"
$ws.Range("C15").Value = "S B MARKET ST"
$ws.Range("E15").Value = 4140
$ws.Range("F15").Value = "T"
$ws.Range("H15").Value = 45249.04213364583
$ws.Range("J15").Value = "10/24/23 12:12"
$ws.Range("K15").Value = "10/24/23 12:12"
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = "$4,140 as of 10/24/2023 10:12:25 AM"
$ws.Range("N15").Value = 4120
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0

# Row 16
$ws.Range("A16").Value = "L474817"
$ws.Range("C16").Value = "SAFETY MARKET"
$ws.Range("E16").Value = 4220
$ws.Range("F16").Value = "T"
$ws.Range("H16").Value = 45237.04213364583
$ws.Range("J16").Value = "10/24/23 14:20"
$ws.Range("K16").Value = "10/24/23 14:20"
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = "$4,240 as of 10/24/2023 10:34:37 AM"
$ws.Range("N16").Value = 3960
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0

# Row 17
$ws.Range("A17").Value = "LK236828"
$ws.Range("C17").Value = "WORLDWIDE AUTOMOTIVE"
$ws.Range("E17").Value = 4300
$ws.Range("F17").Value = "T"
$ws.Range("H17").Value = 45240.04213364583
$ws.Range("J17").Value = "10/23/23 10:42"
$ws.Range("K17").Value = "10/23/23 10:42"
$ws.Range("L17").Value = 20
$ws.Range("M17").Value = "$4,300 as of 10/23/2023 8:42:11 AM"
$ws.Range("N17").Value = 4300
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0

# Row 18
$ws.Range("A18").Value = "L688966"
$ws.Range("C18").Value = "S B WESTERN 108TH MARKET"
$ws.Range("E18").Value = 4680
$ws.Range("F18").Value = "T"
$ws.Range("H18").Value = 45240.04213364583
$ws.Range("J18").Value = "10/24/23 12:38"
$ws.Range("K18").Value = "10/24/23 12:38"
$ws.Range("L18").Value = 80
$ws.Range("M18").Value = "$4,680 as of 10/24/2023 10:38:11 AM"
$ws.Range("N18").Value = 4620
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0

# Row 19
$ws.Range("A19").Value = "L697589"
$ws.Range("C19").Value = "S B DISCOUNT MART"
$ws.Range("E19").Value = 5840
$ws.Range("F19").Value = "T"
$ws.Range("H19").Value = 45232.04213364583
$ws.Range("J19").Value = "10/24/23 14:19"
$ws.Range("K19").Value = "10/24/23 14:19"
$ws.Range("L19").Value = 60
$ws.Range("M19").Value = "$5,900 as of 10/24/2023 10:56:38 AM"
$ws.Range("N19").Value = 5800
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 0

# Row 20
$ws.Range("A20").Value = "L475090"
$ws.Range("C20").Value = "S.B. 2"
$ws.Range("E20").Value = 6820
$ws.Range("F20").Value = "T"
$ws.Range("H20").Value = 45254.04213364583
$ws.Range("J20").Value = "10/24/23 10:37"
$ws.Range("K20").Value = "10/24/23 10:37"
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = "$6,820 as of 10/24/2023 8:37:49 AM"
$ws.Range("N20").Value = 6740
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0

# Row 21
$ws.Range("I21").ClearContents()
$ws.Range("A21").Value = "L474746"
$ws.Range("C21").Value = "ZACATES MARKET"
$ws.Range("E21").Value = 6860
$ws.Range("F21").Value = "T"
$ws.Range("H21").Value = 45271.04213364583
$ws.Range("J21").Value = "10/24/23 15:21"
$ws.Range("K21").Value = "10/24/23 15:18"
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = "$6,900 as of 10/24/2023 10:32:02 AM"
$ws.Range("N21").Value = 6840
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 0

# Row 22
$ws.Range("H22").ClearContents()
$ws.Range("A22").Value = "L682801"
$ws.Range("C22").Value = "SB#5"
$ws.Range("E22").Value = 7840
$ws.Range("F22").Value = "T"
$ws.Range("I22").Value = "ATM Inactive greater than 2000 minutes"
$ws.Range("J22").Value = "09/28/23 15:22"
$ws.Range("K22").Value = "09/28/23 12:14"
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = "$7,840 as of 9/28/2023 12:31:50 PM"
$ws.Range("N22").Value = 7840
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 0

# Row 23
$ws.Range("A23").Value = "L474761"
$ws.Range("C23").Value = "BABS MARKET"
$ws.Range("E23").Value = 7960
$ws.Range("F23").Value = "T"
$ws.Range("H23").Value = 45267.04213364583
$ws.Range("J23").Value = "10/23/23 20:57"
$ws.Range("K23").Value = "10/23/23 20:57"
$ws.Range("L23").Value = 100
$ws.Range("M23").Value = "$7,960 as of 10/23/2023 6:57:34 PM"
$ws.Range("N23").Value = 7900
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 0

# Row 24
$ws.Range("A24").Value = "LK864765"
$ws.Range("C24").Value = "SKY LIQUOR"
$ws.Range("E24").Value = 8740
$ws.Range("F24").Value = "T"
$ws.Range("H24").Value = 45237.04213364583
$ws.Range("J24").Value = "10/24/23 12:53"
$ws.Range("K24").Value = "10/24/23 12:04"
$ws.Range("L24").Value = 60
$ws.Range("M24").Value = "$8,800 as of 10/24/2023 10:04:24 AM"
$ws.Range("N24").Value = 8780
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0

# Row 25
$ws.Range("A25").Value = "L704741"
$ws.Range("C25").Value = "W ADAMS COIN LAUNDRY"
$ws.Range("E25").Value = 9880
$ws.Range("F25").Value = "T"
$ws.Range("H25").Value = 45236.04213364583
$ws.Range("J25").Value = "10/24/23 15:12"
$ws.Range("K25").Value = "10/24/23 15:12"
$ws.Range("L25").Value = 40
$ws.Range("M25").Value = "$9,960 as of 10/24/2023 11:13:30 AM"
$ws.Range("N25").Value = 9660
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 0

# Row 26
$ws.Range("A26").Value = "LK891176"
$ws.Range("C26").Value = "98 DISCOUNT STORE"
$ws.Range("E26").Value = 12000
$ws.Range("F26").Value = "T"
$ws.Range("H26").Value = 45234.04213364583
$ws.Range("J26").Value = "10/24/23 13:58"
$ws.Range("K26").Value = "10/24/23 13:58"
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = "$12,040 as of 10/24/2023 10:03:30 AM"
$ws.Range("N26").Value = 11720
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = 0

# Row 27
$ws.Range("A27").Value = "LK923383"
$ws.Range("C27").Value = "SAMYS PHONE CARDS"
$ws.Range("E27").Value = 17020
$ws.Range("F27").Value = "T"
$ws.Range("H27").Value = 45247.04213364583
$ws.Range("J27").Value = "10/24/23 14:43"
$ws.Range("K27").Value = "10/24/23 14:43"
$ws.Range("L27").Value = 80
$ws.Range("M27").Value = "$17,120 as of 10/24/2023 10:08:28 AM"
$ws.Range("N27").Value = 16920
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 0
$ws.Range("Q27").Value = 0

# Row 28
$ws.Range("A28").Value = "Total Outstanding Cash Balance:"
$ws.Range("E28").Value = 130420
